$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 57884
$ws.Range("A3").Value = 130999073
$ws.Range("B3").Value = 57884
$ws.Range("E3").Value = 100109
$ws.Range("F3").Value = "Tretåig hackspett"
$ws.Range("G3").Value = "Picoides tridactylus"
$ws.Range("H3").Value = "(Linnaeus, 1758)"
$ws.Range("Q3").Value = 460090
$ws.Range("R3").Value = 7058767
$ws.Range("AC3").Value = "Ringhack äldre"
$ws.Range("A4").Value = 130999075
$ws.Range("B4").Value = 79243
$ws.Range("E4").Value = 6425
$ws.Range("F4").Value = "Garnlav"
$ws.Range("G4").Value = "Alectoria sarmentosa"
$ws.Range("H4").Value = "(Ach.) Ach."
$ws.Range("Q4").Value = 459962
$ws.Range("R4").Value = 7058771
$ws.Range("AC4").Value = "Mycket rikligt"
$ws.Range("B5").Value = 57884
$ws.Range("B6").Value = 79243
$ws.Range("B7").Value = 57884
$ws.Range("B8").Value = 57884
$ws.Range("B9").Value = 79243
$ws.Range("B10").Value = 79243
$ws.Range("B11").Value = 79243
$ws.Range("B12").Value = 57884
$ws.Range("B13").Value = 79243
$ws.Range("A14").Value = 130999067
$ws.Range("B14").Value = 57884
$ws.Range("Q14").Value = 459832
$ws.Range("R14").Value = 7058625
$ws.Range("A15").Value = 130999065
$ws.Range("B15").Value = 57884
$ws.Range("Q15").Value = 459953
$ws.Range("R15").Value = 7058639
$ws.Range("A16").Value = 130999072
$ws.Range("B16").Value = 57884
$ws.Range("Q16").Value = 460119
$ws.Range("R16").Value = 7058721
$ws.Range("AC16").Value = "Ringhack äldre"
$ws.Range("A17").Value = 130999061
$ws.Range("B17").Value = 57884
$ws.Range("Q17").Value = 459956
$ws.Range("R17").Value = 7058732
$ws.Range("AC17").Value = "Ringhack"
$ws.Range("B18").Value = 79243
$ws.Range("B19").Value = 57884
$ws.Range("B20").Value = 79243
$ws.Range("B21").Value = 57884
$ws.Range("B22").Value = 57988
$ws.Range("B23").Value = 79243
$ws.Range("B24").Value = 79243
